$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: a lookup that failed to find the establishment name/CNPJ
$ws.Range("A5").Value = "NOME NÃO ENCONTRADO"
$ws.Range("B5").Value = "CNPJ NÃO ENCONTRADO"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = "'2025-09-02"
$ws.Range("D5").ClearFormats()

# Row 6: successfully parsed QR code for nfe-sp
$ws.Range("A6").Value = "RESTAURANTE DOM PEDRO LTDA"
$ws.Range("B6").Value = "03.031.196/0001-70"
$ws.Range("C6").Value = 81.90000000000001
$ws.Range("D6").Value = "'2025-09-02"
$ws.Range("D6").ClearFormats()
